$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.877.83"
$ws.Range("E2").Value = "  +2.88%  "
$ws.Range("D3").Value = "2.423.03"
$ws.Range("E3").Value = "  +2.73%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "554.64"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.02"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.93%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("E8").Value = "  +0.93%  "
$ws.Range("E9").Value = "  +4.29%  "
$ws.Range("E10").Value = "  +3.65%  "
$ws.Range("E11").Value = "  +1.52%  "
$ws.Range("E12").Value = "  -2.17%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "24.69"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +3.36%  "
$ws.Range("D14").Value = "2.856.42"
$ws.Range("E14").Value = "  +2.82%  "
$ws.Range("D15").Value = "59.838.09"
$ws.Range("E15").Value = "  +2.84%  "
$ws.Range("E16").Value = "  +4.52%  "
$ws.Range("D17").Value = "2.447.50"
$ws.Range("E17").Value = "  +4.39%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.36"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +5.69%  "
$ws.Range("E19").Value = "  +3.69%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "334.54"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.60%  "
$ws.Range("E21").Value = "  +1.34%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "64.67"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +2.68%  "
$ws.Range("E24").Value = "  +1.42%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.62"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.48%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("E27").Value = "  -1.01%  "
$ws.Range("D28").Value = "0.0₃0789"
$ws.Range("E28").Value = "  +6.84%  "
$ws.Range("E29").Value = "  +3.24%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "170.67"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.28"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +2.13%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.67"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +1.10%  "
$ws.Range("E33").Value = "  +0.83%  "
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("E35").Value = "  +5.71%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.26"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.10%  "
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.63"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.68%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "40.16"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +2.54%  "
$ws.Range("E40").Value = "  +11.70%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "313.40"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +6.63%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.74"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +2.51%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "142.32"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -3.44%  "
$ws.Range("E44").Value = "  +1.74%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0522"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +3.70%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.25"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.12%  "
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.571"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.42%  "
$ws.Range("B48").Value = "Polygon"
$ws.Range("C48").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.404"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +5.87%  "
$ws.Range("E49").Value = "  +2.91%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "11.04"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.33%  "
$ws.Range("E51").Value = "  +4.74%  "
